$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.484.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "3.177.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'534.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'144.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.60%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'7.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.113"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.23%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "Cardano"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.432"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.75%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "3.728.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "TRON"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.140"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'25.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0000173"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "59.532.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "3.192.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'374.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.528"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.35%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'69.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'8.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "0.0₃0891"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'1.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'22.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.97%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.00%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.27%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'156.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "  +7.85%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'25.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "2.720.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.86%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "  +4.22%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "VeChain"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.0294"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.94%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'39.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.04%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "Mantle"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.724"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "3.223.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.991"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "Cosmos"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'6.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "Stellar"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.100"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.95%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'20.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.91%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.768"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.87%  "
$ws.Range("E51").Style = "Normal"
